$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.147.84'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '4.040.78'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '538.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.06%  '
$ws.Range('D7').Value = '4.037.14'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.698'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('E11').Value = '  -2.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.70'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000328'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').Value = '4.676.13'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = '4.046.69'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('E19').Value = '  -2.00%  '
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').Value = '72.101.09'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '436.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '98.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.53'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +24.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.30'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('E30').Value = '  +1.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +23.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.135'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '50.37'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +18.08%  '
$ws.Range('E35').Value = '  -0.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '682.25'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '67.38'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.460'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.94%  '
$ws.Range('D39').Value = '0.0₃0845'
$ws.Range('E39').Value = '  -1.63%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.149'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.19%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.41'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.40'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +17.48%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0494'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.67'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.11'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.57%  '
$ws.Range('D51').Value = '2.874.44'
$ws.Range('E51').Value = '  +10.60%  '
